# Add attendance for Raahil Madhok's presentation (row 2) on the
# "2025 - Spring" sheet: 7 in-person attendees, 4 via Zoom.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025 - Spring")

$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 4
